$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "NAME"
$ws.Range("C1").Value = "CONTACT"
$ws.Range("D1").Value = "MAIL"

# Data row
$ws.Range("B4").Value = "ROHIT KAWADE"
$ws.Range("D4").Value = "rohitkawade021@gmail.com"

$ws.Range("A1").Value = "ROLL"
$ws.Range("A4").Value = 1160
$ws.Range("C4").Value = 9822907360

# Footer note
$ws.Range("A8").Value = "//EDITED THID FILE ALSO"

# Column widths (values tuned so that, after the runtime's internal
# whole-pixel quantization of ColumnWidth, the exported widths land on
# 16 / ~11.43 / 26 characters as in the target workbook)
$ws.Columns("B").ColumnWidth = 15.1
$ws.Columns("C").ColumnWidth = 10.6
$ws.Columns("D").ColumnWidth = 25.1

# Hyperlink on the e-mail cell (also applies the built-in Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:rohitkawade021@gmail.com") | Out-Null

# Final selection left on A8
$ws.Range("A8").Select() | Out-Null
